$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 12, shifting existing rows 12..45 down to 13..46.
$ws.Rows.Item(12).Insert()

# The newly inserted row 12 is blank; populate it by copying the values
# that used to occupy row 12 (now duplicated at row 13 after the shift),
# then overwrite the columns that actually change for the new entry.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(12, $col).Value2 = $ws.Cells.Item(13, $col).Value2
}

$ws.Cells.Item(12, 4).Value2 = 45274
$ws.Cells.Item(12, 14).Value2 = 20000
$ws.Cells.Item(12, 15).Value2 = 20000
$ws.Cells.Item(12, 16).Value2 = 20000
$ws.Cells.Item(12, 19).Value2 = 2000
